$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns retain text formatting so numeric-looking
# strings (e.g. "0.998", "1.00") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "68.928.96"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").Value = "3.372.67"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "586.42"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").Value = "178.16"
$ws.Range("E6").Value = "  +1.94%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").Value = "0.186"
$ws.Range("E9").Value = "  +4.83%  "
$ws.Range("D10").Value = "0.584"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").Value = "48.32"
$ws.Range("E11").Value = "  +6.62%  "
$ws.Range("D12").Value = "0.0000275"
$ws.Range("E12").Value = "  +2.73%  "
$ws.Range("D13").Value = "696.57"
$ws.Range("E13").Value = "  +5.96%  "
$ws.Range("D14").Value = "3.913.50"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "8.50"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").Value = "68.867.25"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.120"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.357.80"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "17.58"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("D20").Value = "11.31"
$ws.Range("E20").Value = "  +3.37%  "
$ws.Range("D21").Value = "0.900"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("E22").Value = "  +4.07%  "
$ws.Range("D23").Value = "17.06"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "100.48"
$ws.Range("E24").Value = "  +2.19%  "
$ws.Range("D25").Value = "3.92"
$ws.Range("E25").Value = "  +2.08%  "
$ws.Range("D26").Value = "2.72"
$ws.Range("E26").Value = "  +2.38%  "
$ws.Range("D27").Value = "9.60"
$ws.Range("E27").Value = "  +4.09%  "
$ws.Range("D28").Value = "33.25"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").Value = "8.62"
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("D30").Value = "7.02"
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("D31").Value = "11.15"
$ws.Range("E31").Value = "  +2.18%  "
$ws.Range("D32").Value = "553.58"
$ws.Range("E32").Value = "  -2.51%  "
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("D34").Value = "58.39"
$ws.Range("E34").Value = "  +3.98%  "
$ws.Range("D35").Value = "3.724.83"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "3.41"
$ws.Range("E37").Value = "  +4.35%  "
$ws.Range("D38").Value = "0.145"
$ws.Range("E38").Value = "  +11.07%  "
$ws.Range("D39").Value = "34.73"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").Value = "3.20"
$ws.Range("E40").Value = "  +2.85%  "
$ws.Range("D41").Value = "2.64"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").Value = "0.0₃0679"
$ws.Range("E42").Value = "  +2.89%  "
$ws.Range("D43").Value = "0.337"
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("D44").Value = "0.0416"
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("E45").Value = "  -3.71%  "
$ws.Range("D46").Value = "2.66"
$ws.Range("E46").Value = "  +2.66%  "
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("D50").Value = "132.09"
$ws.Range("E50").Value = "  +2.10%  "
$ws.Range("E51").Value = "  -0.36%  "
